# Updated cryptos list on Sun Sep 17 11:45:53 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.811.43'
$ws.Range("E2").Value = '  +0.54%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.643.60'
$ws.Range("E3").Value = '  -0.03%  '
$ws.Range("E4").Value = '  +0.37%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.45'
$ws.Range("E5").Value = '  +0.42%  '
$ws.Range("E6").Value = '  -0.50%  '
$ws.Range("E7").Value = '  +0.31%  '
$ws.Range("E8").Value = '  -0.15%  '
$ws.Range("E9").Value = '  +0.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.23'
$ws.Range("E10").Value = '  -0.34%  '
$ws.Range("E11").Value = '  +0.22%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.646.73'
$ws.Range("E12").Value = '  +0.34%  '
$ws.Range("E13").Value = '  -0.77%  '
$ws.Range("E14").Value = '  -0.39%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '64.74'
$ws.Range("E15").Value = '  -1.14%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '26.819.59'
$ws.Range("E16").Value = '  +0.45%  '
$ws.Range("E17").Value = '  -1.18%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '214.78'
$ws.Range("E18").Value = '  -1.14%  '
$ws.Range("E19").Value = '  +0.25%  '
$ws.Range("E20").Value = '  +0.91%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.47'
$ws.Range("E21").Value = '  +9.51%  '
$ws.Range("E22").Value = '  -0.73%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.34'
$ws.Range("E23").Value = '  -1.78%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '146.60'
$ws.Range("E24").Value = '  +0.47%  '
$ws.Range("E25").Value = '  +0.46%  '
$ws.Range("E26").Value = '  -1.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.18'
$ws.Range("E27").Value = '  -0.08%  '
$ws.Range("E28").Value = '  -0.71%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0507'
$ws.Range("E29").Value = '  -1.90%  '
$ws.Range("E30").Value = '  +0.33%  '
$ws.Range("E31").Value = '  -1.02%  '
$ws.Range("E32").Value = '  -0.80%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.288.54'
$ws.Range("E33").Value = '  +0.97%  '
$ws.Range("E34").Value = '  -0.69%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.44'
$ws.Range("E35").Value = '  +1.40%  '
$ws.Range("E36").Value = '  -1.55%  '
$ws.Range("E37").Value = '  +0.31%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.821'
$ws.Range("E38").Value = '  -1.07%  '
$ws.Range("E39").Value = '  +0.21%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.808'
$ws.Range("E40").Value = '  -0.96%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.23'
$ws.Range("E41").Value = '  -0.66%  '
$ws.Range("E42").Value = '  -2.43%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.785.45'
$ws.Range("E43").Value = '  +0.15%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '61.47'
$ws.Range("E44").Value = '  +2.62%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '91.84'
$ws.Range("E45").Value = '  +0.04%  '
$ws.Range("E46").Value = '  +1.04%  '
$ws.Range("E47").Value = '  -1.03%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0520'
$ws.Range("E48").Value = '  +0.64%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.64'
$ws.Range("E49").Value = '  -1.54%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0970'
$ws.Range("E50").Value = '  -0.16%  '
$ws.Range("E51").Value = '  +0.53%  '

Write-Host "Updated cryptos list"
